$wb = $excel.ActiveWorkbook

# Rename the first sheet to reflect the new compile timestamp
$ws = $wb.Worksheets.Item(1)
$ws.Name = "compiled on 2024-03-08 09.20.56"

# Update indicator keys in column A
$ws.Range("A37").Value = "PA1.C7.20-29.M"
$ws.Range("A38").Value = "PA1.C7.20-29.W"
$ws.Range("A313").Value = "PA9.2.C4."
